$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-07-07"

# Update the label for the July row
$ws.Range("A8").Value = "July (through 07-07)"

# Update July row (row 8) values
$ws.Range("D8").Value = 10
$ws.Range("G8").Value = 22
$ws.Range("H8").Value = 35
$ws.Range("I8").Value = 34

# Update Total row (row 9) values
$ws.Range("D9").Value = 400
$ws.Range("G9").Value = 494
$ws.Range("H9").Value = 795
$ws.Range("I9").Value = 840
